$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellValues = @{
    "F2" = 1.82
    "G2" = 1.84
    "H2" = 4.5
    "I2" = 4.8
    "K2" = 4.4
    "P2" = 2.52
    "S2" = 2.32
    "T2" = 1.6
    "U2" = 2.4
    "V2" = 1.26
    "W2" = 2.18
    "Y2" = 26
    "Z2" = 1000
    "AA2" = 100
    "AB2" = 13
    "AD2" = 19.5
    "AE2" = 48
    "AF2" = 13.5
    "AH2" = 16.5
    "AI2" = 48
    "AJ2" = 19.5
    "AL2" = 26
    "AM2" = 70
    "AN2" = 7.6
    "AO2" = 1000
    "G3" = 1.45
    "H3" = 8.800000000000001
    "I3" = 9
    "N3" = 5.4
    "O3" = 1.22
    "Q3" = 1.66
    "S3" = 2.66
    "Y3" = 32
    "AA3" = 260
    "AB3" = 9.800000000000001
    "AF3" = 9
    "AJ3" = 12.5
    "AN3" = 5.7
    "F4" = 2.22
    "G4" = 2.24
    "H4" = 3.35
    "I4" = 3.4
    "J4" = 3.9
    "V4" = 1.41
    "W4" = 1.81
    "Y4" = 17.5
    "AF4" = 16
    "AJ4" = 28
    "AK4" = 20
    "AO4" = 23
    "I5" = 6.6
    "N5" = 6.6
    "I6" = 2.64
    "N6" = 5.2
    "S6" = 2.46
    "U6" = 2.56
    "V6" = 1.61
    "F7" = 1.52
    "G7" = 1.59
    "I7" = 9.199999999999999
    "Q7" = 1.72
    "R7" = 1.45
    "V7" = 1.13
    "F8" = 2.5
    "G8" = 2.52
    "Z8" = 21
    "F9" = 2.92
    "G9" = 2.94
    "H9" = 2.42
    "I9" = 2.44
    "J9" = 3.95
    "K9" = 4
    "P9" = 2.8
    "R9" = 1.75
    "S9" = 2.3
    "T9" = 1.49
    "V9" = 1.69
    "W9" = 1.51
    "X9" = 27
    "AA9" = 34
    "AB9" = 19.5
    "AE9" = 21
    "AF9" = 25
    "AG9" = 14.5
    "AI9" = 25
    "AJ9" = 46
    "AK9" = 26
    "AM9" = 46
    "AN9" = 15
    "AO9" = 11
    "F10" = 9
    "G10" = 9.199999999999999
    "H10" = 1.41
    "I10" = 1.42
    "J10" = 5.5
    "K10" = 5.6
    "Q10" = 1.64
    "V10" = 3.4
    "W10" = 1.12
    "Y10" = 10
    "Z10" = 8.6
    "AJ10" = 290
    "AM10" = 120
    "I11" = 22
    "P11" = 2.88
    "Q11" = 1.51
    "R11" = 1.75
    "AC11" = 19
    "AD11" = 990
    "AE11" = 490
    "AG11" = 13
    "AN11" = 3.4
    "N12" = 8.6
    "P12" = 3.55
    "S12" = 1.91
    "U12" = 2.3
    "Z12" = 120
    "AF12" = 11.5
    "AI12" = 95
    "AM12" = 90
    "AN12" = 3.35
    "F13" = 5.5
    "G13" = 5.6
    "H13" = 1.71
    "I13" = 1.72
    "K13" = 4.3
    "O13" = 1.29
    "Q13" = 1.86
    "V13" = 2.38
    "AA13" = 16.5
    "AE13" = 17
    "AF13" = 42
    "F15" = 2.34
    "G15" = 2.68
    "H15" = 2.58
    "I15" = 3
    "J15" = 3.75
    "K15" = 4.6
    "L15" = 1.25
    "M15" = 1.03
    "N15" = 4.9
    "O15" = 1.19
    "P15" = 2.36
    "Q15" = 1.57
    "R15" = 1.55
    "S15" = 2.4
    "T15" = 1.53
    "U15" = 2.46
    "V15" = 1.5
    "W15" = 1.59
    "X15" = 25
    "Y15" = 17
    "AB15" = 16.5
    "AE15" = 32
    "AF15" = 21
    "AL15" = 1000
    "AM15" = 1000
    "AO15" = 19.5
    "G16" = 2.82
    "Q16" = 2
    "T16" = 1.76
    "AA16" = 65
}

foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value2 = $cellValues[$key]
}
